$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "TestSheet"
$r = $ws3.Range("A1")
$r.Interior.Pattern = 1
$r.Interior.PatternColor = 0
$r.Interior.Color = 65535
Write-Host "done setting interior"
